$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: SCD0303 -> SCD0019
$ws.Name = "SCD0019"

# Update the TC_ID value in B2: DGS-318 -> SCD0019-001
$ws.Range("B2").Value = "SCD0019-001"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns("B").ColumnWidth = 11.666666666666666

# Move the active selection from O2 to B3
$ws.Range("B3").Select()
